# "Them tt chuong sach TA" - add the English textbook chapter title under
# the Vietnamese chapter title on the title slide's subtitle placeholder,
# widen/recenter that placeholder, and justify both paragraphs.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The subtitle placeholder ("Rectangle 3") holds the chapter title.
$shp = $s.Shapes.Placeholders.Item(2)

# Reposition / resize the placeholder (EMU -> points, 1 pt = 12700 EMU).
# Only Left/Width actually move; Top/Height are unchanged by this edit, so
# they are left alone to avoid float round-trip drift.
$shp.Left = 611560 / 12700
$shp.Width = 8352928 / 12700

$tr = $shp.TextFrame.TextRange

# Split the existing single run "Chương 14. Phân lớp và ứng dụng trong tìm
# kiếm" into two runs, breaking right before "kiếm", by rewriting the text
# of the leading 42-character sub-range in place (keeps the same rPr on
# both pieces, matching the target run split).
$tr.Characters(1, 42).Text = "Chương 14. Phân lớp và ứng dụng trong tìm "

# Justify the first paragraph (was right-aligned).
$tr.Paragraphs(1, 1).ParagraphFormat.Alignment = 4

# Append a new paragraph with the English textbook chapter reference;
# InsertAfter preserves the trailing endParaRPr on the new last paragraph.
[void]$tr.InsertAfter("`rIIR.C13. Text classification and Naive Bayes")

# Justify the new second paragraph as well.
$tr.Paragraphs(2, 1).ParagraphFormat.Alignment = 4
